# Refresh the "cryptos" price list (GitHub Actions scheduled update).
# Column D ("Price") cells hold plain text that sometimes looks numeric
# (e.g. "485.82", "1.00", "0.0973", with '.' used as a thousands grouping
# mark for big prices like "55.920.95"). Setting .Value on such a string
# directly would let Excel auto-convert it into a real number, which both
# loses the literal text formatting and mangles values like "1.956.11".
# So for column D we briefly force the "@" (Text) number format, assign
# the literal string, then restore the cell's original style ("Normal")
# so no stray formatting/style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '55.920.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.497.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.36%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +11.60%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.510'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.512.47'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.08%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.62%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.331'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.927.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '55.886.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.75%  '
$ws.Range('E17').Value = '  +5.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.509.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.46%  '
$ws.Range('E20').Value = '  +9.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '319.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.410'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.166'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.614.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0784'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.24%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.68'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.865'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.14'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.612'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.993'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0554'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.11%  '
$ws.Range('E44').Value = '  +6.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.82'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +12.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '262.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +23.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0226'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0905'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.934.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.94%  '
